$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value (kept as text to match the sheet's
# existing text-formatted Price/Volume columns).
$updates = @(
    @{ Cell = 'D2'; Value = '317.03' }
    @{ Cell = 'E2'; Value = '1.61%' }
    @{ Cell = 'D3'; Value = '37.82' }
    @{ Cell = 'E3'; Value = '0.39%' }
    @{ Cell = 'D4'; Value = '5.179' }
    @{ Cell = 'E4'; Value = '0.98%' }
    @{ Cell = 'D5'; Value = '0.08039' }
    @{ Cell = 'E5'; Value = '1.80%' }
    @{ Cell = 'E6'; Value = '1.84%' }
    @{ Cell = 'D7'; Value = '8.511' }
    @{ Cell = 'E7'; Value = '2.86%' }
    @{ Cell = 'D8'; Value = '1.924' }
    @{ Cell = 'E8'; Value = '1.06%' }
    @{ Cell = 'E9'; Value = '4.96%' }
    @{ Cell = 'D10'; Value = '0.9380' }
    @{ Cell = 'E10'; Value = '2.12%' }
    @{ Cell = 'D11'; Value = '0.1288' }
    @{ Cell = 'E11'; Value = '9.17%' }
    @{ Cell = 'D12'; Value = '0.1939' }
    @{ Cell = 'E12'; Value = '0.44%' }
    @{ Cell = 'D13'; Value = '0.09087' }
    @{ Cell = 'E13'; Value = '0.01%' }
    @{ Cell = 'D14'; Value = '0.03398' }
    @{ Cell = 'E14'; Value = '2.39%' }
    @{ Cell = 'D15'; Value = '0.09544' }
    @{ Cell = 'E15'; Value = '-0.63%' }
    @{ Cell = 'D16'; Value = '0.001402' }
    @{ Cell = 'E16'; Value = '0.56%' }
    @{ Cell = 'D17'; Value = '0.006168' }
    @{ Cell = 'E17'; Value = '5.34%' }
    @{ Cell = 'E18'; Value = '-4.62%' }
    @{ Cell = 'D19'; Value = '0.3521' }
    @{ Cell = 'E19'; Value = '2.30%' }
    @{ Cell = 'D20'; Value = '6.590' }
    @{ Cell = 'E20'; Value = '24.22%' }
    @{ Cell = 'E21'; Value = '2.44%' }
    @{ Cell = 'D22'; Value = '0.2310' }
    @{ Cell = 'E22'; Value = '-10.71%' }
    @{ Cell = 'D23'; Value = '0.04396' }
    @{ Cell = 'E23'; Value = '0.56%' }
    @{ Cell = 'E24'; Value = '-1.38%' }
    @{ Cell = 'D25'; Value = '0.004269' }
    @{ Cell = 'E25'; Value = '-8.62%' }
    @{ Cell = 'E26'; Value = '-2.14%' }
    @{ Cell = 'D27'; Value = '0.0003988' }
    @{ Cell = 'E27'; Value = '0.05%' }
    @{ Cell = 'D39'; Value = '0.02353' }
    @{ Cell = 'E39'; Value = '4.28%' }
    @{ Cell = 'D40'; Value = '0.05155' }
    @{ Cell = 'E40'; Value = '1.05%' }
    @{ Cell = 'D41'; Value = '0.007695' }
    @{ Cell = 'E41'; Value = '3.40%' }
    @{ Cell = 'D42'; Value = '0.1403' }
    @{ Cell = 'E42'; Value = '3.70%' }
    @{ Cell = 'D43'; Value = '0.008710' }
    @{ Cell = 'E43'; Value = '-3.12%' }
    @{ Cell = 'E44'; Value = '5.57%' }
    @{ Cell = 'D45'; Value = '0.008842' }
    @{ Cell = 'E45'; Value = '2.85%' }
    @{ Cell = 'D46'; Value = '0.00006467' }
    @{ Cell = 'E46'; Value = '-1.38%' }
    @{ Cell = 'E47'; Value = '0.05%' }
    @{ Cell = 'E48'; Value = '-5.22%' }
    @{ Cell = 'E49'; Value = '69.03%' }
    @{ Cell = 'E50'; Value = '0.05%' }
    @{ Cell = 'E51'; Value = '0.05%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = $origStyle
}
